$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$modelText = "MultiOutputRegressor(estimator=GridSearchCV(cv=5,`n                                            estimator=Pipeline(steps=[('model',`n                                                                       GradientBoostingRegressor())]),`n                                            param_grid={'model__max_depth': [3,`n                                                                             5,`n                                                                             7],`n                                                        'model__n_estimators': [50,`n                                                                                100,`n                                                                                150]},`n                                            scoring='neg_mean_squared_error'))"

# Header
$ws.Range("F1").Value = "Modelo"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill model text for rows 2-10
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 6).Value = $modelText
}

# Updated numeric values
$ws.Range("B2").Value = 0.440184382832954
$ws.Range("C2").Value = 0.9282894203889445
$ws.Range("D2").Value = 0.5114954490057118

$ws.Range("B3").Value = 2.878893579773739
$ws.Range("C3").Value = 0.9588727023492973
$ws.Range("D3").Value = 1.266131410638986

$ws.Range("B4").Value = 1.092344271455278
$ws.Range("C4").Value = 0.9460426192480635
$ws.Range("D4").Value = 0.80350905318119

$ws.Range("B5").Value = 1.562193423765018
$ws.Range("D5").Value = 0.938958148630534

$ws.Range("B7").Value = 1.271977756333025
$ws.Range("D7").Value = 0.8344729122263083

$ws.Range("B9").Value = 4.894435758873637
$ws.Range("C9").Value = 0.9414904026042551
$ws.Range("D9").Value = 1.680814134582989

$ws.Range("B10").Value = 0.8471066322618547
$ws.Range("C10").Value = 0.9974748013669915
$ws.Range("D10").Value = 0.7278178842185484
